$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format of the existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add new headers for columns I and J
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Fill data rows 2-25: I column is constant 1, J column mirrors column H (IP)
for ($r = 2; $r -le 25; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipValue
}
